# Stream graph fail cause of data / new interactive stream graph
#
# gapUS_06_14.csv (sheet 2) had its "rank/score" column (B) stored as plain
# shared-string text fractions ("23/115", "31/128", ...). Bring it in line
# with gendergap_06_14.csv (sheet 1, column D) and gapCN_06_14.csv
# (sheet 3, column B), which already store this as a real formula
# (numerator/denominator) with the matching numeric/bordered style.

$wb = $excel.ActiveWorkbook

$wsGap = $wb.Worksheets.Item(1)   # gendergap_06_14.csv
$wsUS  = $wb.Worksheets.Item(2)   # gapUS_06_14.csv
$wsCN  = $wb.Worksheets.Item(3)   # gapCN_06_14.csv

# --- fix gapUS_06_14.csv column B: text fractions -> real formulas ---
$wsUS.Range("B2").Formula  = "=23/115"
$wsUS.Range("B3").Formula  = "=31/128"
$wsUS.Range("B4").Formula  = "=27/130"
$wsUS.Range("B5").Formula  = "=31/134"
$wsUS.Range("B6").Formula  = "=19/134"
$wsUS.Range("B7").Formula  = "=17/135"
$wsUS.Range("B8").Formula  = "=22/135"
$wsUS.Range("B9").Formula  = "=23/136"
$wsUS.Range("B10").Formula = "=20/142"

# Match the formatting already used for the equivalent column on the other
# two sheets (numeric, bordered, centered) instead of plain text.
$rng = $wsUS.Range("B2:B10")
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# --- restore each sheet's last on-screen selection ---
$wsGap.Activate()
$wsGap.Range("D2:D10").Select() | Out-Null

$wsUS.Activate()
$wsUS.Range("C41").Select() | Out-Null

$wsCN.Activate()
$wsCN.Range("B1:B1048576").Select() | Out-Null
